$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has one data row (row 2, "Yash"). The new layout has
# 7 data rows: 6 new registrations above it, with the original row pushed
# down to row 8 and its name/timestamp refreshed. Move the existing row's
# values down to row 8 directly (rather than Rows.Insert, which would copy
# the header's direct formatting onto the new rows) so the new rows end up
# as plain, unstyled data cells like the rest of the sheet.
for ($col = 1; $col -le 5; $col++) {
    $ws.Cells.Item(8, $col).Value = $ws.Cells.Item(2, $col).Text
}

$data = @(
    @("Eve",     "05", "ECE",               "eve@gmail.com",     "06-02-2026 19:28:59"),
    @("David",   "04", "CSE-CyberSecurity",  "david@gmail.com",   "06-02-2026 19:28:42"),
    @("Charlie", "03", "CSE-ML",             "charlie@gmail.com", "06-02-2026 19:28:15"),
    @("Bob",     "02", "CSE-DS",             "bob@gmail.com",     "06-02-2026 19:27:52"),
    @("Alice",   "01", "CSE-AIML",           "alice@gmail.com",   "06-02-2026 19:27:35"),
    @("Sai",     "57", "CSE-AI",             "sai@gmail.com",     "06-02-2026 19:27:06")
)

# Roll numbers like "05" / "04" / "57" look numeric, so Excel's normal type
# inference would silently store them as numbers. Force the whole Roll
# column (B2:B7) to text first, write the values, then restore the default
# style so no direct formatting is left on the cells (matches the rest of
# the sheet, which uses no explicit style on data cells).
$rollRange = $ws.Range("B2:B7")
$rollRange.NumberFormat = "@"

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

$rollRange.Style = "Normal"

# The row that moved down to row 8 needs its name and timestamp updated
# (roll "5A2", branch "CSE", and email are unchanged).
$ws.Cells.Item(8, 1).Value = "Yashwanth"
$ws.Cells.Item(8, 5).Value = "06-02-2026 19:26:32"
